$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtTest40mm")

# --- Header values ---
$ws.Range("C2").Value = 557

# New formula in I3 (percent reduction based on test 7 / I column)
$ws.Range("I3").Formula = "=1-I10/C2"

# Footnote marker on J4 ("*")
$ws.Range("J4").Value = "*"

# --- Row 6: muscle length ---
$ws.Range("D6").Value = 9.8042999999999996
$ws.Range("E6").Value = 5.4329000000000001
$ws.Range("F6").Value = 1.5627
$ws.Range("G6").Value = 3.1173000000000002
$ws.Range("H6").Value = 4.5189000000000004
$ws.Range("I6").Value = 1.1431
$ws.Range("J6").Value = 20.042999999999999
$ws.Range("K6").Value = 3.3938000000000001
$ws.Range("L6").Value = 4.0789999999999997
$ws.Range("M6").Value = 20.853000000000002
$ws.Range("N6").Value = 10.808999999999999

# --- Row 7: dl/dtheta ---
$ws.Range("D7").Value = 106
$ws.Range("E7").Value = 96
$ws.Range("F7").Value = 66.5
$ws.Range("G7").Value = 74.5
$ws.Range("H7").Value = 85.5
$ws.Range("I7").Value = 53.5
$ws.Range("J7").Value = 53.5
$ws.Range("K7").Value = 17
$ws.Range("L7").Value = 17
$ws.Range("M7").Value = 42.5
$ws.Range("N7").Value = 30

# --- Row 8: Hand ---
$ws.Range("D8").Value = 29
$ws.Range("E8").Value = 26.5
$ws.Range("F8").Value = 32.299999999999997
$ws.Range("G8").Value = 31.2
$ws.Range("H8").Value = 27.6
$ws.Range("I8").Value = 33.6
$ws.Range("J8").Value = 36.5
$ws.Range("K8").Value = 36
$ws.Range("L8").Value = 35.200000000000003
$ws.Range("M8").Value = 39.9
$ws.Range("N8").Value = 39.700000000000003

# --- Row 9: MA ---
$ws.Range("D9").Value = 34.5
$ws.Range("E9").Value = 31
$ws.Range("F9").Value = 44.5
$ws.Range("G9").Value = 35.5
$ws.Range("H9").Value = 33
$ws.Range("I9").Value = 35.5
$ws.Range("J9").Value = 35
$ws.Range("K9").Value = 35
$ws.Range("L9").Value = 37
$ws.Range("M9").Value = 40.5
$ws.Range("N9").Value = 38

# --- Row 10: ICR (spline) ---
$ws.Range("D10").Value = 529.5
$ws.Range("E10").Value = 520
$ws.Range("F10").Value = 524
$ws.Range("G10").Value = 522
$ws.Range("H10").Value = 522
$ws.Range("I10").Value = 526
$ws.Range("J10").Value = 470
$ws.Range("K10").Value = 461
$ws.Range("L10").Value = 461
$ws.Range("M10").Value = 462
$ws.Range("N10").Value = 457

# --- Row 13: ICR to muscle ---
$ws.Range("D13").Value = 65
$ws.Range("E13").Value = 57
$ws.Range("F13").Value = 73
$ws.Range("G13").Value = 62
$ws.Range("H13").Value = 67
$ws.Range("I13").Value = 74
$ws.Range("J13").Value = 65
$ws.Range("K13").Value = 70
$ws.Range("L13").Value = 75
$ws.Range("M13").Value = 62.5
$ws.Range("N13").Value = 69

# --- New row 16: second footnote marker ---
$ws.Range("J16").Value = "*changed pressure from 100 kPa to 230 kPa"

# --- Cosmetic: selection moved to N14 ---
$ws.Range("N14").Select()

